# This script inserts a new data row at row 95 (pushing all existing
# rows 95-201 down to 96-202), and populates the newly inserted row
# with a new weekly price observation for "Piña" / Feria Lagunitas de
# Puerto Montt, matching the rest of the columns of the row that used
# to occupy row 95 (now row 96).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 95; this shifts rows 95..201
# down to 96..202 and keeps all their values/styles intact.
$ws.Rows.Item(95).Insert()

# Populate the newly-inserted row 95 with the new record.
$ws.Cells.Item(95, 1).Value2  = 4
$ws.Cells.Item(95, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(95, 3).Value2  = "Los Lagos"
$ws.Cells.Item(95, 4).Value2  = 44601
$ws.Cells.Item(95, 5).Value2  = 10
$ws.Cells.Item(95, 6).Value2  = "Fruta"
$ws.Cells.Item(95, 7).Value2  = 100108
$ws.Cells.Item(95, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(95, 9).Value2  = 100108005
$ws.Cells.Item(95, 10).Value2 = "Piña"
$ws.Cells.Item(95, 11).Value2 = "Caramelo"
$ws.Cells.Item(95, 12).Value2 = "Tercera"
$ws.Cells.Item(95, 13).Value2 = 80
$ws.Cells.Item(95, 14).Value2 = 17000
$ws.Cells.Item(95, 15).Value2 = 18000
$ws.Cells.Item(95, 16).Value2 = 17500
$ws.Cells.Item(95, 17).Value2 = "$/caja 16 unidades"
$ws.Cells.Item(95, 18).Value2 = "Ecuador"
$ws.Cells.Item(95, 19).Value2 = 1094
$ws.Cells.Item(95, 20).Value2 = 16
